# Update the "Bill Summary" sheet with the latest billed-quantity figures:
# rows 8-12 (line items) plus the Grand Total / Net Payable rows (14 & 16).
#
# Columns D and G (and H on the total rows) are stored as TEXT even though
# their contents look numeric (e.g. "2", "14592.00") -- assigning a bare
# numeric string via .Value would make Excel auto-convert the cell to a
# Number, so we briefly force a Text number format, assign the literal
# string, then restore the default "Normal" style so the cell's format
# is left exactly as it was (no lingering Text number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 8 - "P. point" / Short point (up to 3 mtr.)
$ws.Range("A8").Value = "P. point"
$ws.Range("C8").Value = 57
Set-TextValue "D8" "2"
$ws.Range("E8").Value = "Short point (up to 3 mtr.)"
$ws.Range("F8").Value = 256
Set-TextValue "G8" "14592.00"

# Row 9 - "P. point" / Medium point (up to 6 mtr.)
$ws.Range("A9").Value = "P. point"
$ws.Range("C9").Value = 10
Set-TextValue "D9" "3"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
Set-TextValue "G9" "4720.00"

# Row 10 - "R. mtr." / 25 mm
$ws.Range("A10").Value = "R. mtr."
$ws.Range("C10").Value = 93
Set-TextValue "D10" "17"
$ws.Range("E10").Value = "25 mm"
$ws.Range("F10").Value = 56
Set-TextValue "G10" "5208.00"

# Row 11 - LED batten
$ws.Range("C11").Value = 66
Set-TextValue "D11" "27"
$ws.Range("E11").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F11").Value = 492
Set-TextValue "G11" "32472.00"

# Row 12 - "%" / Add Tender Premium
$ws.Range("A12").Value = "%"
$ws.Range("C12").Value = 55
Set-TextValue "D12" "37"
$ws.Range("E12").Value = "Add Tender Premium "
$ws.Range("F12").Value = 0
Set-TextValue "G12" "0.00"

# Row 14 - Grand Total Rs.
Set-TextValue "G14" "56992.00"
Set-TextValue "H14" "56992.00"

# Row 16 - NET PAYABLE AMOUNT Rs.
Set-TextValue "G16" "56992.00"
Set-TextValue "H16" "56992.00"
